$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & volume refresh, plus a new
# LidoStakedEther entry that shifted subsequent rows down by one).
# A leading apostrophe forces Excel to store numeric-looking strings as
# plain text (matching the original inlineStr cell type), and resetting
# the Style back to "Normal" afterwards avoids leaving a stray @ text
# number-format style on the cell.

$ws.Range("D2").Value = "'58.182.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.96%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.132.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'533.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'138.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'LidoStakedEther"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'3.132.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.34%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'XRP"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.463"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.39%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'Toncoin"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'7.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.41%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'Cardano"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.410"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.70%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'3.672.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.15%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'TRON"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.137"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'Avalanche"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'25.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.77%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'58.245.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.81%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.136.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'Polkadot"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'6.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.10%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Chainlink"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'12.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'8.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'BitcoinCash"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'355.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'Dai"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'69.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.79%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Polygon"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.505"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Kaspa"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.168"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.82%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Binance-PegBSC-USD"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'PEPE"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0884"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'InternetComputer(DFINITY)"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'7.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.78%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'RenderToken"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'6.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.64%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'1.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.50%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'EthereumClassic"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'21.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.20%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'NEARProtocol"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'5.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.63%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Fetch.AI"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.58%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Monero"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'158.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.67%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Aptos"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'6.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.21%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'EnergySwap"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'25.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.91%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.77%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Stacks"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.33%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Hedera"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.0668"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.701"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Filecoin"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'4.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'OKB"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'37.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.18%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Maker"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.405.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.89%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'RenzoRestakedETH"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.178.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.34%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'FirstDigitalUSD"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'VeChain"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0269"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'ONDO"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.980"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.62%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Cosmos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'6.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'19.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.91%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'SuiNetwork"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.742"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.00%  "
$ws.Range("E51").Style = "Normal"
